# Apply the cryptocurrency price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D (Price) and E (Volume(1h)) columns hold plain text that looks numeric
# (e.g. "57.733.29", "  +2.27%  "). A leading apostrophe forces Excel to
# store the literal text instead of silently coercing it to a number/percent.

$ws.Range("D2").Value = "'" + '57.733.29'
$ws.Range("E2").Value = "'" + '  +2.27%  '
$ws.Range("D3").Value = "'" + '3.032.18'
$ws.Range("E3").Value = "'" + '  +1.93%  '
$ws.Range("E4").Value = "'" + '  +0.00%  '
$ws.Range("D5").Value = "'" + '511.80'
$ws.Range("E5").Value = "'" + '  +1.96%  '
$ws.Range("D6").Value = "'" + '139.84'
$ws.Range("E6").Value = "'" + '  +3.81%  '
$ws.Range("E7").Value = "'" + '  +0.06%  '
$ws.Range("E8").Value = "'" + '  +3.43%  '
$ws.Range("D9").Value = "'" + '7.50'
$ws.Range("E9").Value = "'" + '  +0.84%  '
$ws.Range("D11").Value = "'" + '0.367'
$ws.Range("E11").Value = "'" + '  +5.32%  '
$ws.Range("D12").Value = "'" + '3.550.11'
$ws.Range("E12").Value = "'" + '  +1.93%  '
$ws.Range("E13").Value = "'" + '  +1.94%  '
$ws.Range("D14").Value = "'" + '26.76'
$ws.Range("E14").Value = "'" + '  +5.37%  '
$ws.Range("D15").Value = "'" + '0.0000166'
$ws.Range("E15").Value = "'" + '  +10.36%  '
$ws.Range("D16").Value = "'" + '57.751.95'
$ws.Range("E16").Value = "'" + '  +2.35%  '
$ws.Range("E17").Value = "'" + '  +9.13%  '
$ws.Range("D18").Value = "'" + '3.032.35'
$ws.Range("E18").Value = "'" + '  +1.95%  '
$ws.Range("D19").Value = "'" + '12.91'
$ws.Range("E19").Value = "'" + '  +4.96%  '
$ws.Range("D20").Value = "'" + '8.02'
$ws.Range("E20").Value = "'" + '  +3.83%  '
$ws.Range("D21").Value = "'" + '332.81'
$ws.Range("E21").Value = "'" + '  +3.54%  '
$ws.Range("D22").Value = "'" + '0.998'
$ws.Range("E22").Value = "'" + '  -0.17%  '
$ws.Range("E23").Value = "'" + '  +1.06%  '
$ws.Range("E24").Value = "'" + '  +6.90%  '
$ws.Range("D25").Value = "'" + '64.66'
$ws.Range("E25").Value = "'" + '  +4.83%  '
$ws.Range("E26").Value = "'" + '  +5.25%  '
$ws.Range("E27").Value = "'" + '  +0.36%  '
$ws.Range("D28").Value = "'" + '0.0₃0933'
$ws.Range("E28").Value = "'" + '  +5.42%  '
$ws.Range("D29").Value = "'" + '6.84'
$ws.Range("E29").Value = "'" + '  +7.00%  '
$ws.Range("D30").Value = "'" + '7.48'
$ws.Range("E30").Value = "'" + '  +11.31%  '
$ws.Range("E31").Value = "'" + '  +3.86%  '
$ws.Range("E32").Value = "'" + '  +3.84%  '
$ws.Range("D33").Value = "'" + '20.82'
$ws.Range("E33").Value = "'" + '  +2.57%  '
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").Value = "'" + '4.73'
$ws.Range("E34").Value = "'" + '  +6.81%  '
$ws.Range("B35").Value = 'Monero'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D35").Value = "'" + '155.27'
$ws.Range("E35").Value = "'" + '  -1.87%  '
$ws.Range("D36").Value = "'" + '5.87'
$ws.Range("E36").Value = "'" + '  +6.54%  '
$ws.Range("E37").Value = "'" + '  +2.41%  '
$ws.Range("D38").Value = "'" + '24.87'
$ws.Range("E38").Value = "'" + '  +8.29%  '
$ws.Range("E39").Value = "'" + '  +2.64%  '
$ws.Range("D40").Value = "'" + '3.068.03'
$ws.Range("E40").Value = "'" + '  +1.99%  '
$ws.Range("D41").Value = "'" + '37.40'
$ws.Range("E41").Value = "'" + '  +3.32%  '
$ws.Range("D42").Value = "'" + '3.88'
$ws.Range("E42").Value = "'" + '  +9.51%  '
$ws.Range("E43").Value = "'" + '  +0.05%  '
$ws.Range("D44").Value = "'" + '2.311.84'
$ws.Range("E44").Value = "'" + '  +3.23%  '
$ws.Range("D45").Value = "'" + '0.656'
$ws.Range("E45").Value = "'" + '  +2.78%  '
$ws.Range("E46").Value = "'" + '  +2.80%  '
$ws.Range("D47").Value = "'" + '0.993'
$ws.Range("E47").Value = "'" + '  +1.68%  '
$ws.Range("D48").Value = "'" + '6.03'
$ws.Range("E48").Value = "'" + '  +5.62%  '
$ws.Range("E49").Value = "'" + '  +2.82%  '
$ws.Range("D50").Value = "'" + '19.76'
$ws.Range("E50").Value = "'" + '  +4.68%  '
$ws.Range("E51").Value = "'" + '  -3.57%  '
